# Updates cryptos price/volume data per upstream source refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'34.614.10"
$ws.Range("E2").Value = "'  +0.52%  "

# Row 3
$ws.Range("D3").Value = "'1.818.50"
$ws.Range("E3").Value = "'  +1.47%  "

# Row 4
$ws.Range("E4").Value = "'  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'228.62"
$ws.Range("E5").Value = "'  +1.05%  "

# Row 6
$ws.Range("D6").Value = "'0.580"
$ws.Range("E6").Value = "'  +4.06%  "

# Row 7
$ws.Range("E7").Value = "'  +0.13%  "

# Row 8
$ws.Range("D8").Value = "'35.08"
$ws.Range("E8").Value = "'  +7.41%  "

# Row 9
$ws.Range("D9").Value = "'0.302"
$ws.Range("E9").Value = "'  +1.71%  "

# Row 10
$ws.Range("E10").Value = "'  +0.56%  "

# Row 11
$ws.Range("D11").Value = "'0.0952"
$ws.Range("E11").Value = "'  +0.21%  "

# Row 12
$ws.Range("D12").Value = "'2.081.57"
$ws.Range("E12").Value = "'  +1.48%  "

# Row 13
$ws.Range("E13").Value = "'  +2.95%  "

# Row 14
$ws.Range("D14").Value = "'1.807.17"
$ws.Range("E14").Value = "'  +0.41%  "

# Row 15
$ws.Range("D15").Value = "'0.649"
$ws.Range("E15").Value = "'  +2.06%  "

# Row 16
$ws.Range("D16").Value = "'34.585.77"
$ws.Range("E16").Value = "'  +0.48%  "

# Row 17
$ws.Range("E17").Value = "'  +2.35%  "

# Row 18
$ws.Range("D18").Value = "'69.53"
$ws.Range("E18").Value = "'  +1.03%  "

# Row 19
$ws.Range("E19").Value = "'  +0.25%  "

# Row 20
$ws.Range("D20").Value = "'246.20"
$ws.Range("E20").Value = "'  -0.34%  "

# Row 21
$ws.Range("D21").Value = "'11.58"
$ws.Range("E21").Value = "'  +3.32%  "

# Row 22
$ws.Range("E22").Value = "'  +0.22%  "

# Row 23
$ws.Range("E23").Value = "'  +0.72%  "

# Row 24
$ws.Range("D24").Value = "'172.38"
$ws.Range("E24").Value = "'  +4.74%  "

# Row 25
$ws.Range("D25").Value = "'2.11"
$ws.Range("E25").Value = "'  +2.00%  "

# Row 26
$ws.Range("D26").Value = "'7.56"
$ws.Range("E26").Value = "'  +4.54%  "

# Row 27
$ws.Range("D27").Value = "'16.85"
$ws.Range("E27").Value = "'  +2.01%  "

# Row 28
$ws.Range("E28").Value = "'  +1.99%  "

# Row 29
$ws.Range("E29").Value = "'  -0.02%  "

# Row 30
$ws.Range("D30").Value = "'4.00"
$ws.Range("E30").Value = "'  +2.69%  "

# Row 31
$ws.Range("E31").Value = "'  +1.98%  "

# Row 32
$ws.Range("B32").Value = "'Filecoin"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.85"
$ws.Range("E32").Value = "'  +1.23%  "

# Row 33
$ws.Range("B33").Value = "'PancakeSwap"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.25"
$ws.Range("E33").Value = "'  +1.43%  "

# Row 34
$ws.Range("D34").Value = "'1.85"
$ws.Range("E34").Value = "'  +1.68%  "

# Row 35
$ws.Range("D35").Value = "'1.403.46"
$ws.Range("E35").Value = "'  -1.58%  "

# Row 36
$ws.Range("E36").Value = "'  -0.78%  "

# Row 37
$ws.Range("D37").Value = "'0.681"
$ws.Range("E37").Value = "'  +1.93%  "

# Row 38
$ws.Range("E38").Value = "'  +0.65%  "

# Row 39
$ws.Range("E39").Value = "'  -0.38%  "

# Row 40
$ws.Range("B40").Value = "'Aave"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'83.27"
$ws.Range("E40").Value = "'  -1.77%  "

# Row 41
$ws.Range("B41").Value = "'MXToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.86"
$ws.Range("E41").Value = "'  +4.35%  "

# Row 42
$ws.Range("D42").Value = "'0.956"
$ws.Range("E42").Value = "'  +2.11%  "

# Row 43
$ws.Range("D43").Value = "'2.41"
$ws.Range("E43").Value = "'  +0.18%  "

# Row 44
$ws.Range("D44").Value = "'13.91"
$ws.Range("E44").Value = "'  +2.33%  "

# Row 45
$ws.Range("E45").Value = "'  +2.78%  "

# Row 46
$ws.Range("D46").Value = "'0.0513"
$ws.Range("E46").Value = "'  -1.92%  "

# Row 47
$ws.Range("D47").Value = "'6.05"
$ws.Range("E47").Value = "'  -0.56%  "

# Row 48
$ws.Range("D48").Value = "'1.980.89"
$ws.Range("E48").Value = "'  +1.66%  "

# Row 49
$ws.Range("D49").Value = "'105.62"
$ws.Range("E49").Value = "'  +0.14%  "

# Row 50
$ws.Range("B50").Value = "'PaxDollar"
$ws.Range("C50").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "'  +0.14%  "

# Row 51
$ws.Range("B51").Value = "'BabyDogeCoin"
$ws.Range("C51").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.0₆0130"
$ws.Range("E51").Value = "'  +1.03%  "

Write-Output "Applied all cell updates"